$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Header updates: report volume/number and the week-covering date range
# (rich-text strings collapse to plain text of identical visible content,
# since every run shares the same character formatting).
# ---------------------------------------------------------------------
$ws.Range("A8").Value = "Volume 31   Number  50"
$ws.Range("C9").Value = "Report Covering the Week  12/9/2024  Through  12/15/2024"

# ---------------------------------------------------------------------
# Crime-statistics table updates (rows 15-30)
# ---------------------------------------------------------------------

function Set-NumCell($addr, $value) {
    $ws.Range($addr).Value = $value
}

function Set-TextCell($addr, $text, $styleDonor) {
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = $text
    $ws.Range($styleDonor).Copy() | Out-Null
    $ws.Range($addr).PasteSpecial(-4122) | Out-Null
}

function Set-NumCellWithStyle($addr, $value, $styleDonor) {
    $ws.Range($addr).Value = $value
    $ws.Range($styleDonor).Copy() | Out-Null
    $ws.Range($addr).PasteSpecial(-4122) | Out-Null
}

# Row 15
Set-TextCell "C15" "0" "C14"
Set-NumCell "E15" -100
Set-NumCell "G15" 2
Set-NumCell "H15" 0
Set-NumCell "J15" 11
Set-NumCell "K15" 63.636363636363

# Row 16
Set-NumCell "C16" 4
Set-NumCell "D16" 2
Set-NumCell "E16" 100
Set-NumCell "F16" 12
Set-NumCell "H16" 20
Set-NumCell "I16" 171
Set-NumCell "J16" 161
Set-NumCell "K16" 6.211180124223
Set-NumCell "L16" 0
Set-NumCell "M16" -33.976833976834
Set-NumCell "N16" -83.541867179980

# Row 17
Set-NumCell "C17" 5
Set-NumCell "D17" 5
Set-NumCell "E17" 0
Set-NumCell "F17" 17
Set-NumCell "G17" 26
Set-NumCell "H17" -34.615384615384
Set-NumCell "I17" 344
Set-NumCell "J17" 354
Set-NumCell "K17" -2.824858757062
Set-NumCell "L17" -4.972375690607
Set-NumCell "M17" 67.804878048780
Set-NumCell "N17" -52.155771905424

# Row 18
Set-NumCell "C18" 4
Set-NumCellWithStyle "D18" 5 "I14"
Set-NumCellWithStyle "E18" -20 "K14"
Set-NumCell "G18" 9
Set-NumCell "H18" 11.111111111111
Set-NumCell "I18" 101
Set-NumCell "J18" 104
Set-NumCell "K18" -2.884615384615
Set-NumCell "L18" -35.668789808917
Set-NumCell "M18" -35.256410256410
Set-NumCell "N18" -89.435146443514

# Row 19
Set-NumCell "C19" 3
Set-NumCell "D19" 9
Set-NumCell "E19" -66.666666666666
Set-NumCell "F19" 19
Set-NumCell "G19" 36
Set-NumCell "H19" -47.222222222222
Set-NumCell "I19" 337
Set-NumCell "J19" 454
Set-NumCell "K19" -25.770925110132
Set-NumCell "L19" -33.791748526522
Set-NumCell "M19" -36.891385767790
Set-NumCell "N19" -49.475262368815

# Row 20
Set-NumCell "C20" 3
Set-NumCell "D20" 3
Set-NumCell "E20" 0
Set-NumCell "F20" 14
Set-NumCell "G20" 8
Set-NumCell "H20" 75
Set-NumCell "I20" 114
Set-NumCell "J20" 109
Set-NumCell "K20" 4.587155963302
Set-NumCell "L20" -3.389830508474
Set-NumCell "M20" 4.587155963302
Set-NumCell "N20" -89.920424403183

# Row 21
Set-NumCell "C21" 19
Set-NumCell "D21" 25
Set-NumCell "E21" -24
Set-NumCell "F21" 74
Set-NumCell "G21" 91
Set-NumCell "H21" -18.681318681318
Set-NumCell "I21" 1087
Set-NumCell "J21" 1199
Set-NumCell "K21" -9.341117597998
Set-NumCell "L21" -18.759342301943
Set-NumCell "M21" -15.540015540015
Set-NumCell "N21" -76.204028021015

# Row 22
Set-TextCell "C22" "0" "C14"
Set-NumCell "F22" 2
Set-NumCell "G22" 1
Set-NumCell "H22" 100
Set-NumCell "L22" -3.448275862068
Set-NumCell "M22" -6.666666666666

# Row 23
Set-NumCell "C23" 5
Set-NumCell "D23" 4
Set-NumCell "E23" 25
Set-NumCell "F23" 8
Set-NumCell "G23" 13
Set-NumCell "H23" -38.461538461538
Set-NumCell "I23" 135
Set-NumCell "J23" 160
Set-NumCell "K23" -15.625
Set-NumCell "L23" -21.511627906976
Set-NumCell "M23" 14.406779661016

# Row 24
Set-NumCell "C24" 30
Set-NumCell "D24" 21
Set-NumCell "E24" 42.857142857142
Set-NumCell "F24" 85
Set-NumCell "G24" 80
Set-NumCell "H24" 6.25
Set-NumCell "I24" 967
Set-NumCell "J24" 957
Set-NumCell "K24" 1.044932079414
Set-NumCell "L24" -17.632027257240
Set-NumCell "M24" -15.323992994746

# Row 25
Set-NumCell "C25" 3
Set-NumCell "D25" 5
Set-NumCell "E25" -40
Set-NumCell "F25" 26
Set-NumCell "G25" 26
Set-NumCell "H25" 0
Set-NumCell "I25" 268
Set-NumCell "J25" 285
Set-NumCell "K25" -5.964912280701
Set-NumCell "L25" -31.282051282051

# Row 26
Set-NumCell "C26" 9
Set-NumCell "D26" 7
Set-NumCell "E26" 28.571428571428
Set-NumCell "G26" 41
Set-NumCell "H26" 39.024390243902
Set-NumCell "I26" 755
Set-NumCell "J26" 585
Set-NumCell "K26" 29.059829059829
Set-NumCell "L26" 41.121495327102
Set-NumCell "M26" 48.915187376725

# Row 27
Set-TextCell "C27" "0" "C14"
Set-NumCell "D27" 2
Set-NumCell "E27" -100
Set-NumCell "G27" 3
Set-NumCell "H27" -33.333333333333
Set-NumCell "J27" 19
Set-NumCell "K27" 42.105263157894

# Row 28
Set-TextCell "C28" "0" "C14"
Set-NumCellWithStyle "D28" 1 "I14"
Set-NumCellWithStyle "E28" -100 "K14"
Set-NumCell "F28" 2
Set-NumCell "G28" 1
Set-NumCell "H28" 100
Set-NumCell "I28" 45
Set-NumCell "J28" 51
Set-NumCell "K28" -11.764705882352
Set-NumCell "L28" -16.666666666666

# Row 29
Set-TextCell "D29" "0" "C14"
Set-TextCell "E29" "***.*" "C14"
Set-TextCell "F29" "0" "C14"
Set-NumCell "G29" 3
Set-NumCell "H29" -100

# Row 30
Set-TextCell "D30" "0" "C14"
Set-TextCell "E30" "***.*" "C14"
Set-TextCell "F30" "0" "C14"
Set-NumCell "G30" 2
Set-NumCell "H30" -100
